$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the input values (D6/D7/D8/D9 are formulas and recalc automatically)
$ws.Range("D3").Value = 278728.21000000002
$ws.Range("D5").Value = 75322.320000000007

# Update the active selection / view to D5 (was D3)
$ws.Range("D5").Select()
